# Append new label-translation rows (EN in column A, ZH in column B)
# immediately below the last populated row, matching the formatting of
# the row directly above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @("Battery Output", "电池产量"),
    @("type", "类型"),
    @("Lithium Iron Phosphate", "磷酸铁锂"),
    @("Ternary materials", "三元材料"),
    @("Total", "总量"),
    @("GWh", "百万千瓦时"),
    @("fossil fuel imports", "化石燃料进口"),
    @("year-on-year", "同比变化"),
    @("Coal, Mt", "煤，百万吨"),
    @("Fossil Gas, bcm", "化石燃气，十亿立方米"),
    @("Oil, Mt", "成品油，百万吨"),
    @("12-month moving sum", "12个月移动总和")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

$formatSource = $ws.Range("A" + $lastRow + ":B" + $lastRow)

for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = $pairs[$i][0]
    $ws.Cells.Item($row, 2).Value = $pairs[$i][1]

    $destRange = $ws.Range("A" + $row + ":B" + $row)
    $formatSource.Copy()
    $destRange.PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
